# Updated cryptos list on Mon May 20 07:28:45 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "66.657.04"
$ws.Range("E2").Value = "  -0.75%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.086.87"
$ws.Range("E3").Value = "  -0.81%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  +0.16%  "

# Row 5 (BNB)
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "575.31"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "

# Row 6 (Solana)
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "176.59"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +2.05%  "

# Row 7 (USDC)
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 (LidoStakedEther)
$ws.Range("D8").Value = "3.088.13"
$ws.Range("E8").Value = "  -0.65%  "

# Row 9 (XRP)
$ws.Range("E9").Value = "  -1.33%  "

# Row 10 (Toncoin)
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.36"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.43%  "

# Row 11 (Dogecoin)
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.150"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -2.10%  "

# Row 12 (Cardano)
$ws.Range("E12").Value = "  -2.88%  "

# Row 13 (ShibaInu)
$ws.Range("E13").Value = "  -3.57%  "

# Row 14 (Avalanche)
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "35.82"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.78%  "

# Row 15 (TRON)
$ws.Range("E15").Value = "  -0.59%  "

# Row 16 (WrappedliquidstakedEther2.0)
$ws.Range("D16").Value = "3.602.71"
$ws.Range("E16").Value = "  -0.64%  "

# Row 17 (WrappedBTC)
$ws.Range("D17").Value = "66.669.19"
$ws.Range("E17").Value = "  -0.57%  "

# Row 18 (Polkadot)
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "6.95"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.05%  "

# Row 19 (Chainlink)
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "16.72"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +1.63%  "

# Row 20 (WrappedEther)
$ws.Range("D20").Value = "3.087.18"
$ws.Range("E20").Value = "  -0.65%  "

# Row 21 (BitcoinCash)
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "479.76"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "

# Row 22 (Uniswap)
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "7.72"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.47%  "

# Row 23 (Polygon)
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.687"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -2.54%  "

# Row 24 (Litecoin)
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.34"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.71%  "

# Row 25 (InternetComputer(DFINITY))
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.61"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -4.28%  "

# Row 26 (Fetch.AI)
$ws.Range("E26").Value = "  -3.24%  "

# Row 27 (RenderToken)
$ws.Range("E27").Value = "  -4.00%  "

# Row 28 (Dai)
$ws.Range("E28").Value = "  +0.13%  "

# Row 29 (NEARProtocol)
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.92"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.34%  "

# Row 30 (ImmutableX)
$ws.Range("E30").Value = "  -3.97%  "

# Row 31 (PancakeSwap)
$ws.Range("E31").Value = "  -2.66%  "

# Row 32 (EthereumClassic)
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "27.88"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.77%  "

# Row 33 (Hedera)
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.112"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "

# Row 34 (PEPE)
$ws.Range("D34").Value = "0.0₃0938"
$ws.Range("E34").Value = "  -0.76%  "

# Row 35 (FirstDigitalUSD)
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +0.26%  "

# Row 36 (Arweave)
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "48.25"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +2.56%  "

# Row 37 (Filecoin)
$ws.Range("E37").Value = "  -5.13%  "

# Row 38 (Mantle)
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.937"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -3.86%  "

# Row 39 (OKB)
$ws.Range("E39").Value = "  -2.23%  "

# Row 40 (TheGraph)
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.309"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "

# Row 41 (Stacks(now))
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.99"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.77%  "

# Row 42 (Kaspa(now))
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.122"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.51%  "

# Row 43 (Cosmos)
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "8.28"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.13%  "

# Row 44 (dogwifhat)
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.66"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +2.82%  "

# Row 45 (Maker)
$ws.Range("D45").Value = "2.781.63"
$ws.Range("E45").Value = "  -0.81%  "

# Row 46 (Bittensor)
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "370.03"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -3.52%  "

# Row 47 (Monero)
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "135.30"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.11%  "

# Row 48 (VeChain)
$ws.Range("E48").Value = "  -2.61%  "

# Row 49 (USDe)
$ws.Range("E49").Value = "  -0.01%  "

# Row 50 (InjectiveProtocol)
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "24.62"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "

# Row 51 (ThetaToken)
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.22"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.94%  "
